$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 15:30"

# 2. Row 2 (TankONO) gets a fresh price read:
#    - previous current price (B2=35.9) shifts into "Old Cena" (C2)
#    - new current price (B2) = 36.5
#    - Delta Cena (D2) is written as a formatted text string "+0.6"
#    - Old Datum (E2) is written as a plain text timestamp (not yet
#      normalised back to a numeric Excel date/time by the script)
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = 36.5
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+0.6"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-02-25 15:30:03"

# 3. Row 3 (Tesco): the previously text-formatted Delta/Old Datum values get
#    normalised back into real numeric values (number + Excel date serial),
#    matching the style already used for the other "Old Datum" cells.
$ws.Range("D3").Value = 0.2

$ws.Range("E3").NumberFormat = $ws.Range("E4").NumberFormat
$ws.Range("E3").Value = 44617.63697916667
